$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 250, shifting existing rows 250:306 down to 251:307
$ws.Rows.Item(250).Insert(-4121)

# Populate the newly inserted row 250 with the new data record
$ws.Cells.Item(250, 1).Value = 8
$ws.Cells.Item(250, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(250, 3).Value = "Coquimbo"
$ws.Cells.Item(250, 4).Value = 45204
$ws.Cells.Item(250, 5).Value = 4
$ws.Cells.Item(250, 6).Value = 100112001
$ws.Cells.Item(250, 7).Value = "Berenjena"
$ws.Cells.Item(250, 8).Value = "Sin especificar"
$ws.Cells.Item(250, 9).Value = "Primera"
$ws.Cells.Item(250, 10).Value = 520
$ws.Cells.Item(250, 11).Value = 9000
$ws.Cells.Item(250, 12).Value = 10000
$ws.Cells.Item(250, 13).Value = 9500
$ws.Cells.Item(250, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(250, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(250, 16).Value = 190
$ws.Cells.Item(250, 17).Value = 50
$ws.Cells.Item(250, 18).Value = "Hortaliza"
